$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The fixture had a few stray extra rows (old rows 5-7) left over from earlier
# scratch testing - drop them before reshaping the real table.
$ws.Rows("5:7").Delete()

# Make room for a new leading index column; this shifts the old ID/Password/Type
# columns from A:C one slot to the right, into B:D (along with their formatting).
$ws.Columns("A:A").Insert()

# The old D column (now E) only ever held an empty, formatted placeholder cell -
# drop it now that it no longer has a purpose.
$ws.Range("E1:E4").Clear()

# New column A is a simple 0-based row index for the data rows.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2

# Header row + new index column share the same bold/centered/bordered style.
$headerRange = $ws.Range("B1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexRange = $ws.Range("A2:A4")
$indexRange.Font.Bold = $true
$indexRange.Borders.LineStyle = 1
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160

# Data cells (now in B:D for rows 2-4) go back to plain/general formatting.
$ws.Range("B2:D4").ClearFormats()

# Park the selection back on A1 (the sheet was mid-scroll/selected into the
# rows that just got removed).
$ws.Range("A1").Select()
